$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5655272006988525
$ws.Range("B1").Value = 0.4628998339176178
$ws.Range("C1").Value = 4.969462871551514
$ws.Range("D1").Value = 2.825546264648438
$ws.Range("E1").Value = 1.225590229034424
